$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values would otherwise be auto-converted to numbers
$textCells = @("D5", "D6", "D7", "D9", "D10", "D15", "D20", "D21", "D22", "D24", "D25", "D27", "D30", "D31", "D32", "D33", "D34", "D36", "D38", "D40", "D41", "D42", "D49", "D50", "D51")
foreach ($cellref in $textCells) {
    $ws.Range($cellref).NumberFormat = "@"
}

# Apply updated cell values from the commit diff
$ws.Range('D2').Value = '57.298.89'
$ws.Range('E2').Value = '  -0.59%  '
$ws.Range('D3').Value = '3.096.45'
$ws.Range('E3').Value = '  -0.10%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '524.25'
$ws.Range('E5').Value = '  +0.16%  '
$ws.Range('D6').Value = '136.58'
$ws.Range('E6').Value = '  -3.37%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').Value = '3.097.22'
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').Value = '0.448'
$ws.Range('E9').Value = '  +2.28%  '
$ws.Range('D10').Value = '7.30'
$ws.Range('E10').Value = '  +1.37%  '
$ws.Range('E11').Value = '  -0.95%  '
$ws.Range('E12').Value = '  +2.34%  '
$ws.Range('D13').Value = '3.634.77'
$ws.Range('E13').Value = '  +0.03%  '
$ws.Range('E14').Value = '  +2.55%  '
$ws.Range('D15').Value = '25.23'
$ws.Range('E15').Value = '  -1.62%  '
$ws.Range('E16').Value = '  -0.59%  '
$ws.Range('D17').Value = '57.417.70'
$ws.Range('E17').Value = '  -0.56%  '
$ws.Range('D18').Value = '3.099.88'
$ws.Range('E18').Value = '  +0.03%  '
$ws.Range('E19').Value = '  -2.51%  '
$ws.Range('D20').Value = '12.34'
$ws.Range('E20').Value = '  -3.41%  '
$ws.Range('D21').Value = '7.83'
$ws.Range('E21').Value = '  -2.55%  '
$ws.Range('D22').Value = '345.12'
$ws.Range('E22').Value = '  +1.74%  '
$ws.Range('E23').Value = '  -0.08%  '
$ws.Range('D24').Value = '67.65'
$ws.Range('E24').Value = '  +1.48%  '
$ws.Range('D25').Value = '0.498'
$ws.Range('E25').Value = '  -2.66%  '
$ws.Range('E26').Value = '  -1.74%  '
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  -0.12%  '
$ws.Range('D28').Value = '0.0₃0886'
$ws.Range('E28').Value = '  -3.31%  '
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('D30').Value = '7.36'
$ws.Range('E30').Value = '  +2.60%  '
$ws.Range('D31').Value = '1.86'
$ws.Range('E31').Value = '  -0.16%  '
$ws.Range('D32').Value = '6.02'
$ws.Range('E32').Value = '  -7.09%  '
$ws.Range('D33').Value = '20.76'
$ws.Range('E33').Value = '  -0.82%  '
$ws.Range('D34').Value = '4.91'
$ws.Range('E34').Value = '  +6.47%  '
$ws.Range('E35').Value = '  -4.15%  '
$ws.Range('D36').Value = '158.75'
$ws.Range('E36').Value = '  +1.97%  '
$ws.Range('E37').Value = '  -1.15%  '
$ws.Range('D38').Value = '25.86'
$ws.Range('E38').Value = '  -4.56%  '
$ws.Range('E39').Value = '  -1.65%  '
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').Value = '0.0657'
$ws.Range('E40').Value = '  -0.64%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').Value = '1.60'
$ws.Range('E41').Value = '  +5.27%  '
$ws.Range('D42').Value = '4.10'
$ws.Range('E42').Value = '  +3.94%  '
$ws.Range('E43').Value = '  +2.23%  '
$ws.Range('D44').Value = '3.137.79'
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range('D45').Value = '2.370.68'
$ws.Range('E45').Value = '  +3.32%  '
$ws.Range('E46').Value = '  -0.78%  '
$ws.Range('E47').Value = '  +0.10%  '
$ws.Range('E48').Value = '  +2.80%  '
$ws.Range('D49').Value = '0.970'
$ws.Range('E49').Value = '  -1.41%  '
$ws.Range('D50').Value = '5.95'
$ws.Range('E50').Value = '  -1.14%  '
$ws.Range('D51').Value = '19.66'
$ws.Range('E51').Value = '  -4.25%  '
